$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as plain text (contains non-numeric characters
# such as extra separators, percent signs, or spacing) -- assign directly.
$textValues = @{
    'D2' = '68.402.63'
    'E2' = '  +0.01%  '
    'D3' = '2.647.72'
    'E3' = '  +0.19%  '
    'E4' = '  +0.04%  '
    'E5' = '  -0.01%  '
    'E6' = '  +2.96%  '
    'E7' = '  +0.06%  '
    'E8' = '  -1.02%  '
    'E9' = '  -1.43%  '
    'E11' = '  +0.45%  '
    'E12' = '  -0.94%  '
    'E13' = '  -0.74%  '
    'D14' = '3.132.01'
    'E14' = '  +0.30%  '
    'E15' = '  -2.99%  '
    'D16' = '68.254.76'
    'E16' = '  -0.11%  '
    'D17' = '2.700.24'
    'E17' = '  +2.17%  '
    'E18' = '  -0.36%  '
    'E19' = '  -1.20%  '
    'E20' = '  -1.53%  '
    'E21' = '  +0.76%  '
    'E22' = '  -2.77%  '
    'E23' = '  -0.10%  '
    'E24' = '  +0.02%  '
    'E25' = '  +0.05%  '
    'E26' = '  -0.90%  '
    'D27' = '2.785.58'
    'E27' = '  +0.16%  '
    'E28' = '  -2.45%  '
    'E30' = '  -2.05%  '
    'E31' = '  -1.79%  '
    'E32' = '  -2.34%  '
    'E33' = '  +0.92%  '
    'E34' = '  +2.81%  '
    'E35' = '  +0.03%  '
    'E36' = '  -2.00%  '
    'E37' = '  +1.47%  '
    'E38' = '  -1.46%  '
    'E39' = '  -0.98%  '
    'E40' = '  -1.70%  '
    'E41' = '  -1.73%  '
    'E42' = '  -1.51%  '
    'D43' = '0.0₆0321'
    'E43' = '  -4.96%  '
    'E44' = '  +0.04%  '
    'E45' = '  +0.53%  '
    'E46' = '  +0.60%  '
    'E47' = '  +0.34%  '
    'E48' = '  -1.59%  '
    'E49' = '  -1.99%  '
    'E50' = '  +0.95%  '
    'E51' = '  -0.21%  '
}
foreach ($ref in $textValues.Keys) {
    $ws.Range($ref).Value = $textValues[$ref]
}

# Cells whose new value looks like a plain number (e.g. "27.99") -- Excel would
# otherwise auto-convert these to numeric cells, so force a Text number format
# before assigning, matching the source workbook which stores them as strings.
$numericTextValues = @{
    'D5' = '597.85'
    'D6' = '159.16'
    'D10' = '0.157'
    'D13' = '27.99'
    'D15' = '0.0000187'
    'D18' = '11.38'
    'D19' = '360.14'
    'D20' = '7.40'
    'D24' = '74.64'
    'D28' = '0.0000103'
    'D30' = '562.24'
    'D34' = '1.64'
    'D35' = '0.999'
    'D37' = '19.69'
    'D38' = '158.56'
    'D41' = '5.34'
    'D45' = '157.52'
    'D46' = '3.81'
    'D47' = '22.04'
    'D48' = '1.69'
    'D49' = '0.0773'
    'D50' = '0.575'
}
foreach ($ref in $numericTextValues.Keys) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $numericTextValues[$ref]
}
